$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new values look like plain numbers (e.g. "0.997") must be
# forced to remain TEXT cells (matching the original inlineStr/shared-string
# type in the workbook) instead of being auto-coerced to numeric cells by
# Excel. We do this, per cell, by temporarily applying a text number format,
# assigning the string value, then clearing the format again so the cell
# keeps its original (default/no) style.

$ws.Range("D2").Value = '28.027.31'
$ws.Range("E2").Value = '  +3.24%  '
$ws.Range("D3").Value = '1.572.52'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -1.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.53'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.19'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +5.53%  '
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0881'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("D12").Value = '1.798.03'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '1.586.22'
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = '28.006.36'
$ws.Range("E16").Value = '  +3.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.35'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.10'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +5.86%  '
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.43'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.997'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.94'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.95'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.20'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.57'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.14'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.57%  '
$ws.Range("D34").Value = '1.416.42'
$ws.Range("E34").Value = '  -2.46%  '
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("E36").Value = '  -5.00%  '
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.540'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.43'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.95%  '
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.996'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.66'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.973'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.82%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.81'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.64'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("D47").Value = '1.710.16'
$ws.Range("E47").Value = '  +0.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.89'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.16%  '
$ws.Range("E49").Value = '  +3.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0524'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0941'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.64%  '
